# Generate Report for Handback
# - Row 7 (the 7a80a1b3-e747-44fd-aba8-d93eab920cfa file) failed the handback
#   transform because the handback file name didn't match the handoff file
#   name. Update the Status column (C) on both locale sheets and record the
#   error detail in a new column L.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C7").Value = "Handback transform failed"
$zh.Range("L7").Value = "Handback file name: 1gbm0uv5.p5l is different with handoff file name: 7a80a1b3-e747-44fd-aba8-d93eab920cfa.b109bbd3c2a7eedbd16164aa78afd93f3804e6a2.zh-cn."

$de = $wb.Worksheets.Item("de-de")
$de.Range("C7").Value = "Handback transform failed"
$de.Range("L7").Value = "Handback file name: 1gbm0uv5.p5l is different with handoff file name: 7a80a1b3-e747-44fd-aba8-d93eab920cfa.b109bbd3c2a7eedbd16164aa78afd93f3804e6a2.de-de."
